{"js": "// Update the two-digit \u00f7 one-digit division prompts in the practice table.\n// Every \"NN\u00f7N=\" string in the document is unique, so a body.search() on the\n// exact text reliably targets the single correct cell. The pairs are\n// applied in REVERSE document order so a freshly-written replacement value\n// can never be re-matched by a later search for a different (originally\n// earlier-in-list) pair -- this matters because \"85\u00f77=\" is replaced with\n// \"21\u00f74=\", while a separate, later cell that originally held \"21\u00f74=\" must\n// become \"59\u00f79=\" (not be overwritten a second time).\nconst replacements = [\n  [\"32\u00f78=\", \"97\u00f78=\"],\n  [\"96\u00f76=\", \"37\u00f75=\"],\n  [\"21\u00f74=\", \"59\u00f79=\"],\n  [\"12\u00f75=\", \"71\u00f77=\"],\n  [\"26\u00f79=\", \"44\u00f77=\"],\n  [\"25\u00f75=\", \"69\u00f78=\"],\n  [\"37\u00f77=\", \"93\u00f73=\"],\n  [\"17\u00f79=\", \"63\u00f73=\"],\n  [\"15\u00f73=\", \"36\u00f78=\"],\n  [\"75\u00f72=\", \"78\u00f78=\"],\n  [\"96\u00f79=\", \"19\u00f78=\"],\n  [\"40\u00f79=\", \"86\u00f75=\"],\n  [\"62\u00f75=\", \"44\u00f77=\"],\n  [\"12\u00f78=\", \"42\u00f79=\"],\n  [\"38\u00f78=\", \"42\u00f72=\"],\n  [\"85\u00f77=\", \"21\u00f74=\"],\n  [\"75\u00f75=\", \"79\u00f75=\"],\n  [\"50\u00f75=\", \"17\u00f75=\"],\n  [\"26\u00f75=\", \"59\u00f77=\"],\n  [\"43\u00f79=\", \"62\u00f77=\"],\n  [\"53\u00f79=\", \"65\u00f73=\"],\n  [\"94\u00f72=\", \"87\u00f74=\"],\n  [\"49\u00f73=\", \"85\u00f73=\"],\n  [\"14\u00f73=\", \"64\u00f78=\"],\n  [\"36\u00f76=\", \"27\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit \u00f7 one-digit division prompts in the practice table.\n# Each \"NN\u00f7N=\" string in the document is unique, so Find/Replace on the\n# exact text reliably targets the correct cell. The pairs are applied in\n# REVERSE document order so that a freshly-written replacement value can\n# never be re-matched by a later (originally earlier-in-list) search --\n# this matters because \"85\u00f77=\" is replaced with \"21\u00f74=\", and a different,\n# later cell originally containing \"21\u00f74=\" must become \"59\u00f79=\".\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $wdFindContinue = 1\n    $wdReplaceOne = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne)\n}\n\nReplace-ExactText \"32\u00f78=\" \"97\u00f78=\"\nReplace-ExactText \"96\u00f76=\" \"37\u00f75=\"\nReplace-ExactText \"21\u00f74=\" \"59\u00f79=\"\nReplace-ExactText \"12\u00f75=\" \"71\u00f77=\"\nReplace-ExactText \"26\u00f79=\" \"44\u00f77=\"\nReplace-ExactText \"25\u00f75=\" \"69\u00f78=\"\nReplace-ExactText \"37\u00f77=\" \"93\u00f73=\"\nReplace-ExactText \"17\u00f79=\" \"63\u00f73=\"\nReplace-ExactText \"15\u00f73=\" \"36\u00f78=\"\nReplace-ExactText \"75\u00f72=\" \"78\u00f78=\"\nReplace-ExactText \"96\u00f79=\" \"19\u00f78=\"\nReplace-ExactText \"40\u00f79=\" \"86\u00f75=\"\nReplace-ExactText \"62\u00f75=\" \"44\u00f77=\"\nReplace-ExactText \"12\u00f78=\" \"42\u00f79=\"\nReplace-ExactText \"38\u00f78=\" \"42\u00f72=\"\nReplace-ExactText \"85\u00f77=\" \"21\u00f74=\"\nReplace-ExactText \"75\u00f75=\" \"79\u00f75=\"\nReplace-ExactText \"50\u00f75=\" \"17\u00f75=\"\nReplace-ExactText \"26\u00f75=\" \"59\u00f77=\"\nReplace-ExactText \"43\u00f79=\" \"62\u00f77=\"\nReplace-ExactText \"53\u00f79=\" \"65\u00f73=\"\nReplace-ExactText \"94\u00f72=\" \"87\u00f74=\"\nReplace-ExactText \"49\u00f73=\" \"85\u00f73=\"\nReplace-ExactText \"14\u00f73=\" \"64\u00f78=\"\nReplace-ExactText \"36\u00f76=\" \"27\u00f75=\"\n"}
